$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(177).Insert()

$ws.Range("A177").Value = 4
$ws.Range("B177").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C177").Value = "Los Lagos"
$ws.Range("D177").Value = 44637
$ws.Range("E177").Value = 10
$ws.Range("F177").Value = 100112032
$ws.Range("G177").Value = "Zapallo italiano"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 120
$ws.Range("K177").Value = 18000
$ws.Range("L177").Value = 18000
$ws.Range("M177").Value = 18000
$ws.Range("N177").Value = '$/caja 50 unidades'
$ws.Range("O177").Value = "Región Metropolitana"
$ws.Range("P177").Value = 360
$ws.Range("Q177").Value = 50
$ws.Range("R177").Value = "Hortaliza"
